$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "ScanDate" header column (J1)
$ws.Range("J1").Value = "ScanDate"

# Record the scan timestamp for processed hosts (rows 2, 4, 5).
# Row 3 (NOTPRIME) is a hostname mismatch, so it is skipped.
$scanDate = "2025-09-09 08:23:02"
$ws.Range("J2").Value = $scanDate
$ws.Range("J4").Value = $scanDate
$ws.Range("J5").Value = $scanDate
